$d = $word.ActiveDocument

# Locate the paragraph containing the "fragments may also contain" text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*fragments may also contain*") {
        $target = $p
    }
}

$r1 = $target.Range.InsertParagraphAfter()
$r2 = $r1.InsertParagraphAfter()

# Now find the second newly-inserted (now empty) paragraph and set its text.
$newPara = $d.Paragraphs.Item(6)
$newPara.Range.InsertAfter([string]::Format("{0}include {1}static{2}%{3}", "{%", [char]0x201C, [char]0x201D, "}"))

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output ("[" + $i + "] style=" + $p.Range.ParagraphStyle.NameLocal + " text=" + $p.Range.Text)
}
